$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'leggings for cold',
    'leggings for men winter',
    'leggings for older women',
    'leggings for pregnant',
    'leggings for running women',
    'leggings for skiing',
    'leggings for tall',
    'leggings for teens',
    'leggings for tennis',
    'leggings for the cold',
    'leggings for thick thighs',
    'leggings for winter',
    'leggings for winter women',
    'leggings for womem',
    'leggings for women 7 for 40',
    'leggings for women 7 pack',
    'leggings for women black pattern',
    'leggings for women brown color',
    'leggings for women capri',
    'leggings for women capri high waist',
    'leggings for women capri length',
    'leggings for women design',
    'leggings for women for cold weather',
    'leggings for women full length',
    'leggings for women knee length',
    'leggings for women light grey',
    'leggings for women long length',
    'leggings for women lot',
    'leggings for women low rise',
    'leggings for women lululemon',
    'leggings for women over 50',
    'leggings for women pack',
    'leggings for women pack of 12',
    'leggings for women pack of 4',
    'leggings for women pack of 7',
    'leggings for women plus size capri',
    'leggings for women sheer',
    'leggings for women shiny',
    'leggings for women short',
    'leggings for women size 16',
    'leggings for women sport',
    'leggings for women sport gym',
    'leggings for women tights',
    'leggings for women with',
    'leggings forever 21',
    'leggings free shipping',
    'leggings free to live',
    'leggings from victorias secret',
    'leggings full length',
    'leggings gradient',
    'leggings gym',
    'leggings gym women',
    'leggings high rise',
    'leggings high waist capri',
    'leggings high waist for women',
    'leggings insulated women',
    'leggings knee',
    'leggings knee length',
    'leggings long',
    'leggings long length for women',
    'leggings lot',
    'leggings lululemon',
    'leggings mid rise',
    'leggings mid waist',
    'leggings net',
    'leggings pack',
    'leggings pack for women',
    'leggings pants',
    'leggings polyester spandex',
    'leggings rainbow',
    'leggings reductores',
    'leggings running',
    'leggings short for women',
    'leggings short length',
    'leggings size small women',
    'leggings sport for women',
    'leggings sport for women high weight',
    'leggings sports',
    'leggings sports nike',
    'leggings support',
    'leggings tall',
    'leggings tight',
    'leggings tights',
    'leggings training',
    'leggings training women',
    'leggings underwear women',
    'leggings victoria',
    'leggings victoria secret',
    'leggings victoria secret women',
    'leggings with',
    'leggings with belly support',
    'leggings with copper',
    'leggings with designs for women',
    'leggings with feet',
    'leggings with foot straps',
    'leggings with grip',
    'leggings with holes',
    'leggings with holes down the side',
    'leggings with holes for women',
    'leggings with holes in the knee'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
